# Add "2022-Q4" data: a new per-quarter sheet plus a new summary row on
# "总计" (Total), matching the upstream "feat: add 2022-Q4 data" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (Total) sheet: insert a new row 2 for the 2022-Q4 summary.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Copy formatting from the (now-shifted) old row 2 so the new row matches
# the look of the existing rows (bordered/centered index cell in col A).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 23
$total.Range("D2").Value = 7.32

# Column A is a 0-based running index; renumber it for every data row now
# that a row has been inserted at the top.
for ($r = 2; $r -le 8; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" sheet with the fund holdings table, placed right
#    after "总计" (i.e. before the sheet that used to be "2022-Q3").
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q3"))
$newSheet.Name = "2022-Q4"

# NOTE: re-fetch sheet references by name after Worksheets.Add() - a
# handle obtained before the insert can end up pointing at the wrong
# sheet once the collection shifts.
$newSheet = $wb.Worksheets.Item("2022-Q4")

# Header row: copy style from the 2022-Q3 sheet's header, then set text.
$wb.Worksheets.Item("2022-Q3").Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows: copy style from the 2022-Q3 sheet's first data row (gives
# column A the bordered/centered index style used throughout the book).
$wb.Worksheets.Item("2022-Q3").Range("A2:H2").Copy()
$newSheet.Range("A2:H24").PasteSpecial(-4122)

# Columns B-G hold text (fund codes with leading zeros, and numbers kept
# as formatted strings like "29.17") - force text so values round-trip
# exactly instead of being coerced to numbers.
$newSheet.Range("B2:G24").NumberFormat = "@"

$data = @(
    @("000960", "招商医药健康产业股票", "29.17", "93.56", "6.78", "1.9777", 3),
    @("340008", "兴全有机增长混合", "21.32", "79.09", "4.97", "1.0596", 4),
    @("012186", "招商品质成长混合A", "13.00", "92.99", "6.63", "0.8619", 4),
    @("560080", "汇添富中证中药 ETF", "18.27", "99.54", "3.71", "0.6778", 8),
    @("011373", "招商前沿医疗保健股票A", "9.48", "92.63", "6.07", "0.5754", 4),
    @("009360", "招商创新增长混合A", "6.04", "92.89", "6.36", "0.3841", 4),
    @("012187", "招商品质成长混合C", "5.57", "92.99", "6.63", "0.3693", 4),
    @("217009", "招商核心价值混合", "7.76", "84.62", "2.99", "0.2320", 8),
    @("159647", "鹏华中证中药ETF", "6.29", "97.07", "3.62", "0.2277", 8),
    @("217010", "招商大盘蓝筹混合", "6.93", "86.28", "2.99", "0.2072", 9),
    @("011690", "招商品质发现混合A", "6.07", "83.24", "2.99", "0.1815", 8),
    @("013559", "招商均衡回报混合A", "3.62", "88.42", "2.99", "0.1082", 7),
    @("014841", "东方阿尔法医疗健康混合A", "1.58", "89.34", "6.57", "0.1038", 4),
    @("011374", "招商前沿医疗保健股票C", "1.09", "92.63", "6.07", "0.0662", 4),
    @("014840", "招商裕华混合", "2.01", "75.22", "3.20", "0.0643", 5),
    @("009361", "招商创新增长混合C", "0.88", "92.89", "6.36", "0.0560", 4),
    @("561510", "华泰柏瑞中证中药ETF", "1.29", "96.27", "3.68", "0.0475", 7),
    @("014589", "招商成长先导股票A", "0.77", "87.62", "5.92", "0.0456", 4),
    @("562390", "银华中证中药ETF", "0.79", "97.84", "3.64", "0.0288", 8),
    @("014842", "东方阿尔法医疗健康混合C", "0.43", "89.34", "6.57", "0.0283", 4),
    @("013560", "招商均衡回报混合C", "0.20", "88.42", "2.99", "0.0060", 7),
    @("014590", "招商成长先导股票C", "0.09", "87.62", "5.92", "0.0053", 4),
    @("011691", "招商品质发现混合C", "0.12", "83.24", "2.99", "0.0036", 8)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}
